# Updates cryptocurrency price (column D) and 1h volume change (column E)
# figures on the active worksheet, per the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few "Price" cells are decimal strings with a significant trailing zero
# (e.g. "5.30", "159.40", "0.370"). Plain .Value assignment would let Excel
# auto-detect them as numbers and silently drop the trailing zero, so force
# those specific cells to Text first, then clear the format back off again
# so the cell keeps matching its neighbours (no leftover numberformat).
$textFormatRows = 11,37,39
foreach ($r in $textFormatRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "68.348.00"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.648.60"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "597.01"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "159.17"
$ws.Range("E6").Value = "  +2.89%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.542"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "0.142"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "5.30"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "0.351"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "27.97"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "3.128.56"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").Value = "68.164.92"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "2.650.62"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "359.73"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").Value = "7.41"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").Value = "2.06"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "74.51"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "9.72"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "2.778.95"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "560.76"
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("D31").Value = "8.02"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "1.65"
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("D37").Value = "159.40"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "19.69"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "0.370"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "2.62"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").Value = "0.0₆0319"
$ws.Range("E43").Value = "  -6.55%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "157.48"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").Value = "22.02"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").Value = "0.0773"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "0.575"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  -0.25%  "

foreach ($r in $textFormatRows) {
    $ws.Range("D$r").ClearFormats()
}
